# Generate Report for Handback
# Adds a new "test2.md" row to the Overview sheet and the zh-cn detail
# sheet, mirroring the existing Ping.md rows (new shared strings +
# hyperlinks, one new row on each of the two sheets).

$wb = $excel.ActiveWorkbook

# Hyperlink font color used throughout the workbook for linked file-name
# cells: RGB FF6495ED -> OLE/BGR long value expected by Range.Font.Color.
$linkColor = 15570276

# ---------------------------------------------------------------------
# Sheet "Overview": append row 5 (test2.md / zh-cn status)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A5"),
    "https://github.com/OpenLocalizationTest/oltest/blob/c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1/e2e/test2.md",
    "",
    "",
    "test2.md"
)
$wsOverview.Range("A5").Font.Underline = $true
$wsOverview.Range("A5").Font.Color = $linkColor

$wsOverview.Range("B5").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "zh-cn": append row 5 (full handback-status detail for test2.md)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A5"),
    "https://github.com/OpenLocalizationTest/oltest/blob/c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1c1/e2e/test2.md",
    "",
    "",
    "test2.md"
)
$wsZhCn.Range("A5").Font.Underline = $true
$wsZhCn.Range("A5").Font.Color = $linkColor

$wsZhCn.Range("B5").Value = "Handed back: in sync with en-US"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C5"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2d2d2d2d2d2d2d2d2d2d2d2d2d2d2d2d2d2d2d2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/master/test2.dd770ae7d8d0bb37ce5217d18b66c19f089bd53d.zh-cn.xlf",
    "",
    "",
    "test2.dd770ae7d8d0bb37ce5217d18b66c19f089bd53d.zh-cn.xlf"
)
$wsZhCn.Range("C5").Font.Underline = $true
$wsZhCn.Range("C5").Font.Color = $linkColor

$wsZhCn.Range("D5").Value = "2016-02-23 05:04:54"
$wsZhCn.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E5"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e3e3e3e3e3e3e3e3e3e3e3e3e3e3e3e3e3e3e3e3/e2e/test2.md",
    "",
    "",
    "test2.md"
)
$wsZhCn.Range("E5").Font.Underline = $true
$wsZhCn.Range("E5").Font.Color = $linkColor

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f4f4f4f4f4f4f4f4f4f4f4f4f4f4f4f4f4f4f4f4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/master/test2.dd770ae7d8d0bb37ce5217d18b66c19f089bd53d.zh-cn.xlf",
    "",
    "",
    "test2.dd770ae7d8d0bb37ce5217d18b66c19f089bd53d.zh-cn.xlf"
)
$wsZhCn.Range("F5").Font.Underline = $true
$wsZhCn.Range("F5").Font.Color = $linkColor

$wsZhCn.Range("G5").Value = "2016-02-23 05:34:31"
$wsZhCn.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H5").Value = "Include"
